$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# New merger/fusion rows (K-39 .. Mn-55 chain) matching the commit
# "adicionei algumas fusoes e uma lista com os elementos faltantes em game"

$ws.Range('A62').Value = 'Ar-36'
$ws.Range('B62').Value = 'He-4'
$ws.Range('C62').Value = 'K-39, n'
$ws.Range('E62').Value = ' O argônio-36 pode capturar uma partícula alfa, resultando na formação de potássio-39 e liberando um nêutron.'

$ws.Range('A63').Value = 'Ar-38'
$ws.Range('B63').Value = 'H-1'
$ws.Range('C63').Value = 'K-39, y'
$ws.Range('E63').Value = 'O argônio-38 pode capturar um próton para formar potássio-39, emitindo um raio gama.'

$ws.Range('A64').Value = 'K-39'
$ws.Range('B64').Value = 'n'
$ws.Range('C64').Value = 'K-40'
$ws.Range('E64').Value = 'O potássio-39 pode capturar um nêutron livre durante processos de captura lenta de nêutrons (s-process), formando potássio-40.'

$ws.Range('A65').Value = 'K-40'
$ws.Range('B65').Value = 'n'
$ws.Range('C65').Value = 'K-41'
$ws.Range('E65').Value = 'O potássio-40 captura outro nêutron, resultando em potássio-41'

$ws.Range('A66').Value = 'Cr-50'
$ws.Range('B66').Value = 'H-1'
$ws.Range('C66').Value = 'Mn-51, y'

$ws.Range('A67').Value = 'Mn-51'
$ws.Range('C67').Value = 'Fe-51, e+, ve'

$ws.Range('A68').Value = 'Fe-51'
$ws.Range('B68').Value = 'H-1'
$ws.Range('C68').Value = 'Co-52, y'

$ws.Range('A69').Value = 'Co-52'
$ws.Range('C69').Value = 'V-51, e+, ve'
$ws.Range('E69').Value = 'Através de uma série de capturas de prótons e decaimentos beta, cromo-50 pode eventualmente levar à formação de vanádio-51'

$ws.Range('A70').Value = 'Ti-49'
$ws.Range('B70').Value = 'n'
$ws.Range('C70').Value = 'Ti-50'

$ws.Range('A71').Value = 'Ti-50'
$ws.Range('B71').Value = 'H-1'
$ws.Range('C71').Value = 'V-51, y'

$ws.Range('A72').Value = 'V-51'
$ws.Range('B72').Value = 'n'
$ws.Range('C72').Value = 'V-52'

$ws.Range('A73').Value = 'V-52'
$ws.Range('C73').Value = 'Cr-52, e+, ve'

$ws.Range('A74').Value = 'V-51'
$ws.Range('B74').Value = 'e'
$ws.Range('C74').Value = 'V-50, ve'
$ws.Range('E74').Value = 'O vanádio-50 pode ser formado a partir do vanádio-51 através de captura eletrônica durante processos de captura lenta de nêutrons.'

$ws.Range('A75').Value = 'Cr-54'
$ws.Range('B75').Value = 'H-1'
$ws.Range('C75').Value = 'Mn-55, y'
$ws.Range('E75').Value = 'O cromo-54 captura um próton, resultando em manganês-55.'

$ws.Range('A76').Value = 'Mn-54'
$ws.Range('B76').Value = 'n'
$ws.Range('C76').Value = 'Mn-55'
$ws.Range('E76').Value = 'O manganês-54 captura um nêutron para formar manganês-55.'

# View-state tweaks (zoom + selection) matching the saved sheet view
$win = $excel.ActiveWindow
$win.Zoom = 85
$win.ScrollRow = 49
$win.ScrollColumn = 1
$ws.Range('D57').Select()
